$d = $word.ActiveDocument

# The document ends with a paragraph "Dia 04/09: 3hr (1 dia)" right before
# the sectPr. We append a new paragraph after it, matching the same
# paragraph/run formatting (inherited automatically by InsertParagraphAfter),
# containing the new metric line "Dia 05/09: 1hr (1 dia)".

$lastPara = $d.Paragraphs.Last
$endRange = $lastPara.Range
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Dia 05/09: 1hr (1 dia)"

Write-Host "Paragraphs: $($d.Paragraphs.Count)"
Write-Host "New last paragraph: $($d.Paragraphs.Last.Range.Text)"
